$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 19:35"

# Row 4 (Estados Unidos)
$ws.Range("E4").Value = 1032924
$ws.Range("G4").Value = 750
$ws.Range("H4").Value = 84175

# Row 9 (row index 13 -> Peru)
$ws.Range("B9").Value = 180737
$ws.Range("C9").Value = 3135
$ws.Range("E9").Value = 95505
$ws.Range("G9").Value = 231
$ws.Range("H9").Value = 12635

# Row 11 (row index 15)
$ws.Range("B11").Value = 173824
$ws.Range("C11").Value = 653
$ws.Range("E11").Value = 17332
$ws.Range("G11").Value = 54
$ws.Range("H11").Value = 7792

# Row 31 (row index 35)
$ws.Range("B31").Value = 23401
$ws.Range("C31").Value = 159
$ws.Range("D31").Value = 19470
$ws.Range("E31").Value = 2434
$ws.Range("G31").Value = 9
$ws.Range("H31").Value = 1497

# Row 35 (row index 39)
$ws.Range("B35").Value = 16548
$ws.Range("C35").Value = 19
$ws.Range("D35").Value = 12232
$ws.Range("E35").Value = 4052
$ws.Range("G35").Value = 4
$ws.Range("H35").Value = 264
